$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 43: survey_unit / character / description
# (combination of survey with quarter or season)
$ws.Range("A43").Value = "survey_unit"
$ws.Range("C43").Value = "character"
$ws.Range("D43").Value = "combination of survey with quarter or season (useful for BITS, NS-IBTS, SWC-IBTS, NEUS, SEUS, SCS, GMEX)"

# Match the row height used by the rest of the table (15.75pt, custom)
$ws.Rows.Item(43).RowHeight = 15.75

# Scroll the view down so row 28 is at the top and select D44,
# matching the author's on-screen state after adding the row.
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D44").Select()
